$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '37.548.58'
$ws.Range("E2").Value = '  -0.45%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.069.38'
$ws.Range("E3").Value = '  -0.26%  '
$ws.Range("E4").Value = '  +0.14%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '232.42'
$ws.Range("E5").Value = '  -0.49%  '
$ws.Range("E6").Value = '  +1.16%  '
$ws.Range("E7").Value = '  +0.05%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '57.55'
$ws.Range("E8").Value = '  -1.79%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.389'
$ws.Range("E9").Value = '  -1.41%  '
$ws.Range("E10").Value = '  -0.38%  '
$ws.Range("E11").Value = '  +1.57%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '14.93'
$ws.Range("E12").Value = '  +0.93%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '2.376.81'
$ws.Range("E13").Value = '  -0.10%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '20.91'
$ws.Range("E14").Value = '  -0.16%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.761'
$ws.Range("E15").Value = '  -1.52%  '
$ws.Range("E16").Value = '  -0.80%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '2.058.76'
$ws.Range("E17").Value = '  -1.76%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '37.511.57'
$ws.Range("E18").Value = '  -0.34%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '70.53'
$ws.Range("E19").Value = '  -1.26%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '5.95'
$ws.Range("E20").Value = '  -2.96%  '
$ws.Range("E21").Value = '  -0.59%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '228.23'
$ws.Range("E22").Value = '  -0.04%  '
$ws.Range("E23").Value = '  +0.10%  '
$ws.Range("E24").Value = '  +0.65%  '
$ws.Range("E25").Value = '  -1.42%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '9.63'
$ws.Range("E26").Value = '  +6.29%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '169.38'
$ws.Range("E27").Value = '  -1.12%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.132'
$ws.Range("E28").Value = '  -2.97%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '19.42'
$ws.Range("E29").Value = '  -0.32%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.37'
$ws.Range("E30").Value = '  -1.52%  '
$ws.Range("E31").Value = '  +1.21%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '4.60'
$ws.Range("E32").Value = '  -1.70%  '
$ws.Range("E33").Value = '  +0.26%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '4.62'
$ws.Range("E34").Value = '  -0.77%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '2.47'
$ws.Range("E35").Value = '  -0.82%  '
$ws.Range("E36").Value = '  -0.21%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '3.32'
$ws.Range("E37").Value = '  -2.40%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.998'
$ws.Range("E38").Value = '  -0.13%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '5.27'
$ws.Range("E39").Value = '  -1.49%  '
$ws.Range("E40").Value = '  +6.98%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '99.54'
$ws.Range("E41").Value = '  -0.04%  '
$ws.Range("E42").Value = '  +4.89%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.0957'
$ws.Range("E43").Value = '  -1.80%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '2.90'
$ws.Range("E44").Value = '  +0.93%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '1.475.47'
$ws.Range("E45").Value = '  +2.79%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '16.65'
$ws.Range("E46").Value = '  -0.34%  '
$ws.Range("E47").Value = '  -1.64%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '4.02'
$ws.Range("E48").Value = '  -4.37%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '7.25'
$ws.Range("E49").Value = '  -1.88%  '
$ws.Range("E50").Value = '  -1.98%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '2.260.16'
$ws.Range("E51").Value = '  -0.19%  '
